$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report record (Pomelo, Feria Lagunitas de Puerto Montt) was
# logged for 2022-03-09. It belongs chronologically right above the existing
# row 116, so insert a fresh row there (pushing rows 116:244 down to 117:245)
# and populate it with the new record's data.
$ws.Rows.Item(116).Insert()

$row = 116
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44629
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100102
$ws.Cells.Item($row, 8).Value = "Cítricos"
$ws.Cells.Item($row, 9).Value = 100102006
$ws.Cells.Item($row, 10).Value = "Pomelo"
$ws.Cells.Item($row, 11).Value = "Start Ruby"
$ws.Cells.Item($row, 12).Value = "Especial"
$ws.Cells.Item($row, 13).Value = 80
$ws.Cells.Item($row, 14).Value = 14000
$ws.Cells.Item($row, 15).Value = 14000
$ws.Cells.Item($row, 16).Value = 14000
$ws.Cells.Item($row, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 1000
$ws.Cells.Item($row, 20).Value = 14
